# Append one new data row (row 38) to the Adafruit IO export sheet,
# mirroring the existing rows' structure/content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 38

# The existing rows store every value (including numeric-looking ones like
# "25") as plain text. Mark the new row's cells as Text *before* writing the
# values so Excel doesn't auto-convert the numeric-looking "25" into a
# number, then restore the default "Normal" style so no stray per-cell
# formatting is left behind in the saved file.
$rowRange = $ws.Range("A" + $row + ":F" + $row)
$rowRange.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"
$ws.Cells.Item($row, 3).Value = "25"
$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"

$rowRange.Style = "Normal"
